# Updated cryptos list on Sat Nov 18 14:23:09 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and swaps the Filecoin / ImmutableX rows (29 / 30 rank values)
# including their name, link, price and volume cells.
#
# Column D values are plain text in this sheet (e.g. "36.503.68" uses dots
# as thousands separators, not a decimal number), so each D assignment is
# prefixed with a literal leading apostrophe - exactly what typing '242.49
# into a cell does in Excel - to keep Excel from reinterpreting them as
# numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '36.503.68'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = "'" + '1.940.51'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = "'" + '242.49'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = "'" + '0.611'
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = "'" + '56.93'
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("D9").Value = "'" + '0.362'
$ws.Range("E9").Value = '  -1.80%  '
$ws.Range("D10").Value = "'" + '0.0805'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("D12").Value = "'" + '2.222.69'
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").Value = "'" + '21.41'
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("D14").Value = "'" + '0.805'
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("D15").Value = "'" + '13.29'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("D16").Value = "'" + '5.15'
$ws.Range("E16").Value = '  -2.67%  '
$ws.Range("D17").Value = "'" + '1.946.75'
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").Value = "'" + '36.417.09'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = "'" + '69.14'
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("D20").Value = "'" + '0.0₃0853'
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("D21").Value = "'" + '227.01'
$ws.Range("E21").Value = '  -0.82%  '
$ws.Range("D22").Value = "'" + '4.94'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = "'" + '2.38'
$ws.Range("E24").Value = '  -4.65%  '
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").Value = "'" + '9.14'
$ws.Range("E26").Value = '  -3.27%  '
$ws.Range("D27").Value = "'" + '160.33'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("D28").Value = "'" + '0.133'
$ws.Range("E28").Value = '  +12.63%  '
$ws.Range("D29").Value = "'" + '19.15'
$ws.Range("E29").Value = '  -1.83%  '
$ws.Range("D30").Value = "'" + '0.118'
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = "'" + '4.66'
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'" + '1.08'
$ws.Range("E32").Value = '  -5.01%  '
$ws.Range("D33").Value = "'" + '0.0612'
$ws.Range("E33").Value = '  -3.13%  '
$ws.Range("D34").Value = "'" + '4.15'
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("E35").Value = '  -0.30%  '
$ws.Range("D36").Value = "'" + '6.10'
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("E37").Value = '  -0.96%  '
$ws.Range("D38").Value = "'" + '2.19'
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("D39").Value = "'" + '3.25'
$ws.Range("E39").Value = '  +13.00%  '
$ws.Range("D40").Value = "'" + '0.0987'
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("E43").Value = '  -2.98%  '
$ws.Range("D44").Value = "'" + '15.68'
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").Value = "'" + '1.337.49'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("D47").Value = "'" + '85.84'
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("D48").Value = "'" + '7.10'
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").Value = "'" + '2.114.30'
$ws.Range("E50").Value = '  -0.60%  '
$ws.Range("D51").Value = "'" + '42.92'
$ws.Range("E51").Value = '  -5.16%  '
